$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 52.138213
$ws.Cells.Item(2, 8).Value = 156.414639
$ws.Cells.Item(2, 9).Value = 0.2220849502516424
$ws.Cells.Item(2, 10).Value = 0.2220849502516423
$ws.Cells.Item(2, 13).Value = 17.10933733333333
$ws.Cells.Item(2, 14).Value = 51.328012
$ws.Cells.Item(2, 15).Value = 0.3554368716515803
$ws.Cells.Item(2, 16).Value = 0.3554368716515803
$ws.Cells.Item(2, 17).Value = 892.0502741741853
$ws.Cells.Item(2, 18).Value = 8028.452467567668
$ws.Cells.Item(2, 19).Value = 0.07893717995834061
$ws.Cells.Item(2, 20).Value = 0.07893717995834061
$ws.Cells.Item(3, 7).Value = 52.138213
$ws.Cells.Item(3, 8).Value = 156.414639
$ws.Cells.Item(3, 9).Value = 0.2220849502516424
$ws.Cells.Item(3, 10).Value = 0.2220849502516423
$ws.Cells.Item(3, 15).Value = 0.2270123898818874
$ws.Cells.Item(3, 16).Value = 0.2270123898818874
$ws.Cells.Item(3, 17).Value = 569.7396100019224
$ws.Cells.Item(3, 18).Value = 5127.656490017302
$ws.Cells.Item(3, 19).Value = 0.05041603531342542
$ws.Cells.Item(3, 20).Value = 0.0504160353134254
$ws.Cells.Item(4, 7).Value = 52.138213
$ws.Cells.Item(4, 8).Value = 156.414639
$ws.Cells.Item(4, 9).Value = 0.2220849502516424
$ws.Cells.Item(4, 10).Value = 0.2220849502516423
$ws.Cells.Item(4, 13).Value = 11.616679
$ws.Cells.Item(4, 14).Value = 34.850037
$ws.Cells.Item(4, 15).Value = 0.2413299803667016
$ws.Cells.Item(4, 16).Value = 0.2413299803667016
$ws.Cells.Item(4, 17).Value = 605.6728840546269
$ws.Cells.Item(4, 18).Value = 5451.055956491643
$ws.Cells.Item(4, 19).Value = 0.05359575668396876
$ws.Cells.Item(4, 20).Value = 0.05359575668396876
$ws.Cells.Item(5, 7).Value = 52.138213
$ws.Cells.Item(5, 8).Value = 156.414639
$ws.Cells.Item(5, 9).Value = 0.2220849502516424
$ws.Cells.Item(5, 10).Value = 0.2220849502516423
$ws.Cells.Item(5, 13).Value = 8.482576333333332
$ws.Cells.Item(5, 14).Value = 25.447729
$ws.Cells.Item(5, 15).Value = 0.1762207580998305
$ws.Cells.Item(5, 16).Value = 0.1762207580998305
$ws.Cells.Item(5, 17).Value = 442.2663716560923
$ws.Cells.Item(5, 18).Value = 3980.397344904831
$ws.Cells.Item(5, 19).Value = 0.03913597829590756
$ws.Cells.Item(5, 20).Value = 0.03913597829590756
$ws.Cells.Item(6, 9).Value = 0.2867694600645705
$ws.Cells.Item(6, 10).Value = 0.2867694600645705
$ws.Cells.Item(6, 13).Value = 17.10933733333333
$ws.Cells.Item(6, 14).Value = 51.328012
$ws.Cells.Item(6, 15).Value = 0.3554368716515803
$ws.Cells.Item(6, 16).Value = 0.3554368716515803
$ws.Cells.Item(6, 17).Value = 1151.869026629333
$ws.Cells.Item(6, 18).Value = 10366.821239664
$ws.Cells.Item(6, 19).Value = 0.1019284397705637
$ws.Cells.Item(6, 20).Value = 0.1019284397705637
$ws.Cells.Item(7, 9).Value = 0.2867694600645705
$ws.Cells.Item(7, 10).Value = 0.2867694600645705
$ws.Cells.Item(7, 15).Value = 0.2270123898818874
$ws.Cells.Item(7, 16).Value = 0.2270123898818874
$ws.Cells.Item(7, 18).Value = 6621.138809148001
$ws.Cells.Item(7, 19).Value = 0.06510022047439662
$ws.Cells.Item(7, 20).Value = 0.06510022047439662
$ws.Cells.Item(8, 9).Value = 0.2867694600645705
$ws.Cells.Item(8, 10).Value = 0.2867694600645705
$ws.Cells.Item(8, 13).Value = 11.616679
$ws.Cells.Item(8, 14).Value = 34.850037
$ws.Cells.Item(8, 15).Value = 0.2413299803667016
$ws.Cells.Item(8, 16).Value = 0.2413299803667016
$ws.Cells.Item(8, 17).Value = 782.081296996
$ws.Cells.Item(8, 18).Value = 7038.731672964001
$ws.Cells.Item(8, 19).Value = 0.06920606816715243
$ws.Cells.Item(8, 20).Value = 0.06920606816715243
$ws.Cells.Item(9, 9).Value = 0.2867694600645705
$ws.Cells.Item(9, 10).Value = 0.2867694600645705
$ws.Cells.Item(9, 13).Value = 8.482576333333332
$ws.Cells.Item(9, 14).Value = 25.447729
$ws.Cells.Item(9, 15).Value = 0.1762207580998305
$ws.Cells.Item(9, 16).Value = 0.1762207580998305
$ws.Cells.Item(9, 17).Value = 571.0809690653333
$ws.Cells.Item(9, 18).Value = 5139.728721588
$ws.Cells.Item(9, 19).Value = 0.05053473165245768
$ws.Cells.Item(9, 20).Value = 0.05053473165245769
$ws.Cells.Item(10, 7).Value = 64.99978900000001
$ws.Cells.Item(10, 8).Value = 194.999367
$ws.Cells.Item(10, 9).Value = 0.2768693838132169
$ws.Cells.Item(10, 10).Value = 0.2768693838132169
$ws.Cells.Item(10, 13).Value = 17.10933733333333
$ws.Cells.Item(10, 14).Value = 51.328012
$ws.Cells.Item(10, 15).Value = 0.3554368716515803
$ws.Cells.Item(10, 16).Value = 0.3554368716515803
$ws.Cells.Item(10, 17).Value = 1112.103316596489
$ws.Cells.Item(10, 18).Value = 10008.9298493684
$ws.Cells.Item(10, 19).Value = 0.0984095876386705
$ws.Cells.Item(10, 20).Value = 0.0984095876386705
$ws.Cells.Item(11, 7).Value = 64.99978900000001
$ws.Cells.Item(11, 8).Value = 194.999367
$ws.Cells.Item(11, 9).Value = 0.2768693838132169
$ws.Cells.Item(11, 10).Value = 0.2768693838132169
$ws.Cells.Item(11, 15).Value = 0.2270123898818874
$ws.Cells.Item(11, 16).Value = 0.2270123898818874
$ws.Cells.Item(11, 17).Value = 710.2843059670505
$ws.Cells.Item(11, 18).Value = 6392.558753703453
$ws.Cells.Item(11, 19).Value = 0.06285278050456392
$ws.Cells.Item(11, 20).Value = 0.06285278050456392
$ws.Cells.Item(12, 7).Value = 64.99978900000001
$ws.Cells.Item(12, 8).Value = 194.999367
$ws.Cells.Item(12, 9).Value = 0.2768693838132169
$ws.Cells.Item(12, 10).Value = 0.2768693838132169
$ws.Cells.Item(12, 13).Value = 11.616679
$ws.Cells.Item(12, 14).Value = 34.850037
$ws.Cells.Item(12, 15).Value = 0.2413299803667016
$ws.Cells.Item(12, 16).Value = 0.2413299803667016
$ws.Cells.Item(12, 17).Value = 755.081683880731
$ws.Cells.Item(12, 18).Value = 6795.73515492658
$ws.Cells.Item(12, 19).Value = 0.0668168829597844
$ws.Cells.Item(12, 20).Value = 0.0668168829597844
$ws.Cells.Item(13, 7).Value = 64.99978900000001
$ws.Cells.Item(13, 8).Value = 194.999367
$ws.Cells.Item(13, 9).Value = 0.2768693838132169
$ws.Cells.Item(13, 10).Value = 0.2768693838132169
$ws.Cells.Item(13, 13).Value = 8.482576333333332
$ws.Cells.Item(13, 14).Value = 25.447729
$ws.Cells.Item(13, 15).Value = 0.1762207580998305
$ws.Cells.Item(13, 16).Value = 0.1762207580998305
$ws.Cells.Item(13, 17).Value = 551.3656718430603
$ws.Cells.Item(13, 18).Value = 4962.291046587543
$ws.Cells.Item(13, 19).Value = 0.04879013271019802
$ws.Cells.Item(13, 20).Value = 0.04879013271019803
$ws.Cells.Item(14, 7).Value = 50.30497766666667
$ws.Cells.Item(14, 8).Value = 150.914933
$ws.Cells.Item(14, 9).Value = 0.2142762058705703
$ws.Cells.Item(14, 10).Value = 0.2142762058705703
$ws.Cells.Item(14, 13).Value = 17.10933733333333
$ws.Cells.Item(14, 14).Value = 51.328012
$ws.Cells.Item(14, 15).Value = 0.3554368716515803
$ws.Cells.Item(14, 16).Value = 0.3554368716515803
$ws.Cells.Item(14, 17).Value = 860.6848324447996
$ws.Cells.Item(14, 18).Value = 7746.163492003196
$ws.Cells.Item(14, 19).Value = 0.07616166428400552
$ws.Cells.Item(14, 20).Value = 0.07616166428400552
$ws.Cells.Item(15, 7).Value = 50.30497766666667
$ws.Cells.Item(15, 8).Value = 150.914933
$ws.Cells.Item(15, 9).Value = 0.2142762058705703
$ws.Cells.Item(15, 10).Value = 0.2142762058705703
$ws.Cells.Item(15, 15).Value = 0.2270123898818874
$ws.Cells.Item(15, 16).Value = 0.2270123898818874
$ws.Cells.Item(15, 17).Value = 549.7069559511386
$ws.Cells.Item(15, 18).Value = 4947.362603560247
$ws.Cells.Item(15, 19).Value = 0.04864335358950149
$ws.Cells.Item(15, 20).Value = 0.04864335358950148
$ws.Cells.Item(16, 7).Value = 50.30497766666667
$ws.Cells.Item(16, 8).Value = 150.914933
$ws.Cells.Item(16, 9).Value = 0.2142762058705703
$ws.Cells.Item(16, 10).Value = 0.2142762058705703
$ws.Cells.Item(16, 13).Value = 11.616679
$ws.Cells.Item(16, 14).Value = 34.850037
$ws.Cells.Item(16, 15).Value = 0.2413299803667016
$ws.Cells.Item(16, 16).Value = 0.2413299803667016
$ws.Cells.Item(16, 17).Value = 584.3767776558357
$ws.Cells.Item(16, 18).Value = 5259.390998902521
$ws.Cells.Item(16, 19).Value = 0.05171127255579606
$ws.Cells.Item(16, 20).Value = 0.05171127255579606
$ws.Cells.Item(17, 7).Value = 50.30497766666667
$ws.Cells.Item(17, 8).Value = 150.914933
$ws.Cells.Item(17, 9).Value = 0.2142762058705703
$ws.Cells.Item(17, 10).Value = 0.2142762058705703
$ws.Cells.Item(17, 13).Value = 8.482576333333332
$ws.Cells.Item(17, 14).Value = 25.447729
$ws.Cells.Item(17, 15).Value = 0.1762207580998305
$ws.Cells.Item(17, 16).Value = 0.1762207580998305
$ws.Cells.Item(17, 17).Value = 426.7158130041285
$ws.Cells.Item(17, 18).Value = 3840.442317037156
$ws.Cells.Item(17, 19).Value = 0.03775991544126726
$ws.Cells.Item(17, 20).Value = 0.03775991544126726
